$d = $word.ActiveDocument

# 1. Update the letter date: September 19, 2025 -> September 21, 2025
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing address line into two paragraphs:
#    "3020 Lamory Pl, Santa Clara CA 95051" ->
#      "3020 Lamory Pl"
#      "Santa Clara, CA 95051"
$rng = $d.Content
$found = $rng.Find.Execute("3020 Lamory Pl, Santa Clara CA 95051", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "3020 Lamory Pl"
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(1, 1) | Out-Null
    $rng.InsertAfter("Santa Clara, CA 95051")
}

# 3. Remove the now-superfluous blank paragraph right after "Board of Directors"
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Board of Directors", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if ($found2) {
    $para = $rng2.Paragraphs(1)
    $nextPara = $para.Next()
    if ($nextPara.Range.Text -eq "`r" -or $nextPara.Range.Text.Trim() -eq "") {
        $nextPara.Range.Delete()
    }
}
